# summer 24 week 14 inputs
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = 1.09
$ws.Range("G4").Value = 1

$ws.Range("C5").Value = 1.38

$ws.Range("D6").Value = 1.53
$ws.Range("G6").Value = 0.96

$ws.Range("D7").Value = 1.72
